$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.034.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.754.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.20"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.30"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.111"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.86"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -12.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.387"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.158"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.245.39"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.88"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.932.79"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000153"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.759.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.12"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.90"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "360.57"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.554"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.27"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.51"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0935"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.25"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.77"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.39"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.23"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.18"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "331.44"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.24"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.89"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0596"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.80"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0257"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.636"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.28"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.38%  "
